$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z3").Select()
